# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" rows (B16:J22) get re-sorted so the periods run in
# ascending order (2305 .. 2311) instead of descending (2311 .. 2305).
# Every column except "Periodo Mora" (E) and "Valor Mora" (F) already holds
# the same value on every one of those rows, so the only cells whose
# displayed content actually changes are E16:E22 and F16:F22 - this script
# writes the post-sort values directly into those two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending "Periodo Mora" values for rows 16-22 (was 2311..2305 descending).
$periodos = @("2305", "2306", "2307", "2308", "2309", "2310", "2311")

# Matching "Valor Mora" values for the same rows (the single odd value
# 43307 moves from the first row to the last row once sorted ascending).
$valores = @(46400, 46400, 46400, 46400, 46400, 46400, 43307)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
